$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 was missing its status value ("IN_STOCK", same as row 2). Fill it in,
# copying the italic status-column formatting used by the other status cells.
$ws.Range("C3").Value = "IN_STOCK"
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)

# Row 4 was missing its employee name. Add the new employee.
$ws.Range("D4").Value = "Сотрудник 3"

# Update the active selection stored with the sheet to D3:D4 (D3 active).
$ws.Range("D3:D4").Select()

$wb.Save()
